$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: insert this week's new record ahead of the existing
# history (row 56), shifting the rest of the Apio series down by one row.
$ws.Rows.Item(56).Insert()

$ws.Cells.Item(56, 1).Value = 11
$ws.Cells.Item(56, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(56, 3).Value = "Bíobío"
$ws.Cells.Item(56, 4).Value = 44540
$ws.Cells.Item(56, 5).Value = 8
$ws.Cells.Item(56, 6).Value = 100112017
$ws.Cells.Item(56, 7).Value = "Apio"
$ws.Cells.Item(56, 8).Value = "Americana (o)"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 250
$ws.Cells.Item(56, 11).Value = 6500
$ws.Cells.Item(56, 12).Value = 7000
$ws.Cells.Item(56, 13).Value = 6700
$ws.Cells.Item(56, 14).Value = "$/docena de matas"
$ws.Cells.Item(56, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(56, 16).Value = 1117
$ws.Cells.Item(56, 17).Value = 6
$ws.Cells.Item(56, 18).Value = "Hortaliza"
